$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52 - GF# match found to be an address-suffix edge case ("20190468b"),
# so it is recorded as text rather than falling back to "did not close".
# Escrow officer is now resolved for this match.
$ws.Range("K52").Value = "20190468b"
$ws.Range("N52").Value = "Tiffany Minnish"

# Row 53 - GF# now resolved to the matching number (was "did not close");
# escrow officer resolved.
$ws.Range("K53").Value = 20190661
$ws.Range("N53").Value = "Tiffany Minnish"

# Row 54
$ws.Range("K54").Value = 20190864
$ws.Range("N54").Value = "Sherry Dixon"

# Row 55
$ws.Range("K55").Value = 20191017
$ws.Range("N55").Value = "Phyllis Alexander"

# Row 56
$ws.Range("K56").Value = 20191129
$ws.Range("N56").Value = "Sherry Dixon"

# Row 60 - GF# correction (edge-case address match picked a different GF#),
# the "Likely Closed" manual flag no longer applies, and the escrow officer
# changes accordingly.
$ws.Range("K60").Value = 20190470
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = "Phyllis Alexander"

# Row 62
$ws.Range("K62").Value = 20190260
$ws.Range("N62").Value = "Sherry Dixon"

# Row 63
$ws.Range("K63").Value = 20190470
$ws.Range("N63").Value = "Phyllis Alexander"

# Row 64
$ws.Range("K64").Value = 20190470
$ws.Range("N64").Value = "Phyllis Alexander"

# Row 86
$ws.Range("K86").Value = 20191028
$ws.Range("N86").Value = "Kristy Kyle"

# Row 90 - GF# resolved, now flagged as Likely Closed, escrow officer resolved.
$ws.Range("K90").Value = 20191129
$ws.Range("M90").Value = $true
$ws.Range("N90").Value = "Sherry Dixon"
